$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column at the front (becomes column A "TabName") and a new
#    row at the bottom (becomes row 4, the new "FilesTab" entry).
# ---------------------------------------------------------------------------
$ws.Columns("A").Insert()
$ws.Rows("4").Insert()

# ---------------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# ---------------------------------------------------------------------------
# 3. Cypher query text (single-quoted here-strings keep backticks / $ / quotes
#    completely literal).
# ---------------------------------------------------------------------------
$statQuery = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['American Staffordshire Terrier']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$casesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['American Staffordshire Terrier']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`
'@

$samplesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) WHERE demo.breed IN ['American Staffordshire Terrier'] WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

$filesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['American Staffordshire Terrier']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$neo4jFileName = "TC02_Canine_Filter_Breed-AmerStaffd_Neo4jData.xlsx"
$webFileName = "TC02_Canine_Filter_Breed-AmerStaffd_WebData.xlsx"

# ---------------------------------------------------------------------------
# 4. Row 2 - CasesTab
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $casesTabQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = $neo4jFileName
$ws.Range("E2").Value = $webFileName

# ---------------------------------------------------------------------------
# 5. Row 3 - SamplesTab
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $samplesTabQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jFileName
$ws.Range("E3").Value = $webFileName

# ---------------------------------------------------------------------------
# 6. Row 4 - FilesTab (new row)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $filesTabQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = $neo4jFileName
$ws.Range("E4").Value = $webFileName

# ---------------------------------------------------------------------------
# 7. Styling - wrap text on B2:C4 (and header C1), matching the diff's s="1"
# ---------------------------------------------------------------------------
$ws.Range("C1").WrapText = $true
$ws.Range("B2:C4").WrapText = $true

# ---------------------------------------------------------------------------
# 8. Column widths
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 10.90625
$ws.Columns("C").ColumnWidth = 123.36328125

# ---------------------------------------------------------------------------
# 9. Row heights (rows 2 & 3 already carry their original explicit heights
#    through the insert operations above - only the brand new row 4 needs one)
# ---------------------------------------------------------------------------
$ws.Rows("4").RowHeight = 246.5

# ---------------------------------------------------------------------------
# 10. Sheet view - selection & zoom (note: this runtime does not persist
#     ActiveWindow.ScrollRow/TopLeftCell to the saved OOXML, so the
#     "topLeftCell=A3" scroll position from the diff cannot be reproduced
#     here - only selection + zoom are reachable through this object model).
# ---------------------------------------------------------------------------
$ws.Range("A4").Select()
$excel.ActiveWindow.Zoom = 70
